$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new SVR parameter columns (K, L, M) with headers and values,
# mirroring the existing row 1 header / row 2 value layout.
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

$ws.Range("K2").Value = 50
$ws.Range("L2").Value = 0.05
$ws.Range("M2").Value = 20

# Remove the now-unused blank, styled cell at A13 (collapses the empty row).
$ws.Range("A13").Clear()

# Restore the active selection to J9, matching the saved workbook state.
$ws.Range("J9").Select() | Out-Null
